# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# --- Hoja1!A1: update the "Conversión del día" note with the new rates ---
$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.54 = 26013.07 pesos`n✅ 26013.07 pesos = 6.51 = 960.55 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newNote

# --- tasas sheet: update the tasas rate table (rows 10 and 12) ---
$ws2.Range("N10").Value = 153
$ws2.Range("O10").Value = 3980
$ws2.Range("N12").Value = 3997.2
$ws2.Range("O12").Value = 147.6
